# correção de varias coisas fodas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "KOALA PARTS"
$ws.Cells.Item(2, 2).Value = "Fonte Carregador De Bateria Jfa 120a Lit"
$ws.Cells.Item(2, 3).Value = "JFA ELETRONICOS"
$ws.Cells.Item(2, 4).Value = "Acessórios para Veículos / Som Automotivo / Módulos Amplificadores"
$ws.Cells.Item(2, 5).Value = "Não"
$ws.Cells.Item(2, 6).Value = "Clássico"
$ws.Cells.Item(2, 7).Value = 5
$ws.Cells.Item(2, 8).Value = "453,51"
$ws.Cells.Item(2, 9).Value = "2.267,55"

# Row 3
$ws.Cells.Item(3, 1).Value = "SUPERTRIO SOM"
$ws.Cells.Item(3, 2).Value = "Fonte Automotiva 70 Amperes Jfa Storm Carregador Com Cor Preto"
$ws.Cells.Item(3, 3).Value = "JFA ELETRONICOS"
$ws.Cells.Item(3, 4).Value = "Acessórios para Veículos / Som Automotivo / Módulos Amplificadores"
$ws.Cells.Item(3, 5).Value = "Não"
$ws.Cells.Item(3, 6).Value = "Premium"
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = "442,82"
$ws.Cells.Item(3, 9).Value = "442,82"

# Row 4
$ws.Cells.Item(4, 1).Value = "TAMANDARÉ"
$ws.Cells.Item(4, 2).Value = "Voltímetro Sequenciador Automotivo Digit"
$ws.Cells.Item(4, 3).Value = "JFA ELETRONICOS"
$ws.Cells.Item(4, 4).Value = "Acessórios para Veículos / Som Automotivo / Outros"
$ws.Cells.Item(4, 5).Value = "Não"
$ws.Cells.Item(4, 6).Value = "Clássico"
$ws.Cells.Item(4, 7).Value = 7
$ws.Cells.Item(4, 8).Value = "62,90"
$ws.Cells.Item(4, 9).Value = "440,30"
